# Applies the "data up to 25th aug" update to fb-surveyInfluenza (Sheet1).
# - Adds 8 new survey dates (18-25 Aug 2020) as rows 201-208, column A.
# - Backfills a handful of previously-blank cells in rows 164-167 (col J)
#   and rows 194-200 (scattered columns) that arrived with this data refresh.
# - Populates the full state-by-state data for the 8 new date rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date labels (column A) for the 8 added rows ---
$ws.Range("A201").Value = "18 08 2020"
$ws.Range("A202").Value = "19 08 2020"
$ws.Range("A203").Value = "20 08 2020"
$ws.Range("A204").Value = "21 08 2020"
$ws.Range("A205").Value = "22 08 2020"
$ws.Range("A206").Value = "23 08 2020"
$ws.Range("A207").Value = "24 08 2020"
$ws.Range("A208").Value = "25 08 2020"

# Row 164
$ws.Range("J164").Value = 0.8002906

# Row 165
$ws.Range("J165").Value = 0

# Row 166
$ws.Range("J166").Value = 0

# Row 167
$ws.Range("J167").Value = 0

# Row 194
$ws.Range("BE194").Value = 0.2834878

# Row 195
$ws.Range("BE195").Value = 0.4254478

# Row 196
$ws.Range("B196").Value = 0.30587
$ws.Range("O196").Value = 0.3770833
$ws.Range("AG196").Value = 0.5991625
$ws.Range("AU196").Value = 0.5598788
$ws.Range("BA196").Value = 0
$ws.Range("BE196").Value = 0.9364174

# Row 197
$ws.Range("B197").Value = 1.3730995
$ws.Range("K197").Value = 0.3045643
$ws.Range("O197").Value = 1.3400551
$ws.Range("AE197").Value = 0.8363382
$ws.Range("AG197").Value = 0.4589349
$ws.Range("AS197").Value = 1.1090937
$ws.Range("AU197").Value = 0.3163754
$ws.Range("BA197").Value = 0
$ws.Range("BE197").Value = 0.6403233

# Row 198
$ws.Range("B198").Value = 0.7339215
$ws.Range("D198").Value = 1.3411649
$ws.Range("I198").Value = 0.033489
$ws.Range("K198").Value = 0.0863913
$ws.Range("O198").Value = 0
$ws.Range("P198").Value = 1.2222825
$ws.Range("Q198").Value = 1.8553284
$ws.Range("T198").Value = 0.4759215
$ws.Range("Y198").Value = 0.3608313
$ws.Range("AD198").Value = 1.107725
$ws.Range("AE198").Value = 0.222036
$ws.Range("AG198").Value = 0.2999368
$ws.Range("AH198").Value = 0.4358372
$ws.Range("AI198").Value = 0.4581942
$ws.Range("AK198").Value = 0.2782701
$ws.Range("AS198").Value = 0.570408
$ws.Range("AU198").Value = 1.6150628
$ws.Range("BA198").Value = 0.2162672
$ws.Range("BD198").Value = 0.7997014
$ws.Range("BE198").Value = 1.7555594

# Row 199
$ws.Range("B199").Value = 0.72778
$ws.Range("C199").Value = 1.4184935
$ws.Range("D199").Value = 0.4474113
$ws.Range("F199").Value = 0.9744556
$ws.Range("G199").Value = 0.7445027
$ws.Range("H199").Value = 0.3636224
$ws.Range("I199").Value = 0.065863
$ws.Range("K199").Value = 0.5077802
$ws.Range("L199").Value = 0.7948132
$ws.Range("M199").Value = 1.0951921
$ws.Range("O199").Value = 0.4781485
$ws.Range("P199").Value = 1.6592639
$ws.Range("Q199").Value = 0.8013213
$ws.Range("R199").Value = 0.3977598
$ws.Range("S199").Value = 1.166078
$ws.Range("T199").Value = 0.751517
$ws.Range("U199").Value = 1.2695013
$ws.Range("V199").Value = 0.7424414
$ws.Range("W199").Value = 0.1347738
$ws.Range("X199").Value = 0.4648749
$ws.Range("Y199").Value = 0.1138349
$ws.Range("Z199").Value = 0.4302361
$ws.Range("AA199").Value = 0.3111848
$ws.Range("AB199").Value = 0.9173894
$ws.Range("AD199").Value = 2.1405613
$ws.Range("AE199").Value = 0.7270262
$ws.Range("AF199").Value = 0.8709389
$ws.Range("AG199").Value = 1.5909839
$ws.Range("AH199").Value = 1.9303104
$ws.Range("AI199").Value = 0.1689025
$ws.Range("AJ199").Value = 0.2601222
$ws.Range("AK199").Value = 0.4587555
$ws.Range("AL199").Value = 0.5012378
$ws.Range("AM199").Value = 0.4697692
$ws.Range("AN199").Value = 0.7040376
$ws.Range("AO199").Value = 0.857339
$ws.Range("AP199").Value = 0.5545387000000001
$ws.Range("AQ199").Value = 0.3303123
$ws.Range("AS199").Value = 0.0770034
$ws.Range("AT199").Value = 0.3489034
$ws.Range("AU199").Value = 0.8697944
$ws.Range("AV199").Value = 0.5390148
$ws.Range("AW199").Value = 0.6981861
$ws.Range("AX199").Value = 0.525369
$ws.Range("AY199").Value = 0.3925474
$ws.Range("BA199").Value = 0
$ws.Range("BB199").Value = 0.4152205
$ws.Range("BC199").Value = 0.8879707999999999
$ws.Range("BD199").Value = 0.6459195
$ws.Range("BE199").Value = 0

# Row 200
$ws.Range("B200").Value = 1.2779974
$ws.Range("C200").Value = 1.2577199
$ws.Range("D200").Value = 1.0835312
$ws.Range("F200").Value = 0.8700725
$ws.Range("G200").Value = 0.5626233
$ws.Range("H200").Value = 0.5373504
$ws.Range("I200").Value = 0.381454
$ws.Range("K200").Value = 0
$ws.Range("L200").Value = 0.834475
$ws.Range("M200").Value = 1.2814193
$ws.Range("O200").Value = 0.3065523
$ws.Range("P200").Value = 1.4724085
$ws.Range("Q200").Value = 1.4887204
$ws.Range("R200").Value = 0.6808582
$ws.Range("S200").Value = 1.1888372
$ws.Range("T200").Value = 1.9380073
$ws.Range("U200").Value = 0.7908345
$ws.Range("V200").Value = 0.7889442
$ws.Range("W200").Value = 0.4522883
$ws.Range("X200").Value = 0.3159411
$ws.Range("Y200").Value = 0.2864132
$ws.Range("Z200").Value = 0.4125149
$ws.Range("AA200").Value = 0.2622421
$ws.Range("AB200").Value = 0.4103554
$ws.Range("AD200").Value = 0.3965797
$ws.Range("AE200").Value = 0.5296465
$ws.Range("AF200").Value = 0.7433696
$ws.Range("AG200").Value = 0
$ws.Range("AH200").Value = 0.4361341
$ws.Range("AI200").Value = 0.2131843
$ws.Range("AJ200").Value = 0.3731063
$ws.Range("AK200").Value = 0.1966041
$ws.Range("AL200").Value = 0.671884
$ws.Range("AM200").Value = 0.4544926
$ws.Range("AN200").Value = 0.4348961
$ws.Range("AO200").Value = 0.5082346
$ws.Range("AP200").Value = 0.8625952
$ws.Range("AQ200").Value = 0.421924
$ws.Range("AS200").Value = 0.4197301
$ws.Range("AT200").Value = 1.0900262
$ws.Range("AU200").Value = 2.1699582
$ws.Range("AV200").Value = 0.630404
$ws.Range("AW200").Value = 0.7919783
$ws.Range("AX200").Value = 0.7244432
$ws.Range("AY200").Value = 0.3873971
$ws.Range("BA200").Value = 0.5522979
$ws.Range("BB200").Value = 0.2796682
$ws.Range("BC200").Value = 0.4251836
$ws.Range("BD200").Value = 0.8042688

# Row 201
$ws.Range("B201").Value = 0.4118234
$ws.Range("C201").Value = 0.7548515
$ws.Range("D201").Value = 0.8182007
$ws.Range("F201").Value = 0.8063285
$ws.Range("G201").Value = 0.4480248
$ws.Range("H201").Value = 0.5712698
$ws.Range("I201").Value = 0.6023397
$ws.Range("K201").Value = 0.649605
$ws.Range("L201").Value = 0.4162498
$ws.Range("M201").Value = 0.8848845
$ws.Range("O201").Value = 0.2195853
$ws.Range("P201").Value = 0.819151
$ws.Range("Q201").Value = 1.5136949
$ws.Range("R201").Value = 0.3146358
$ws.Range("S201").Value = 0.9431047
$ws.Range("T201").Value = 1.0795621
$ws.Range("U201").Value = 0.5865825
$ws.Range("V201").Value = 1.0513243
$ws.Range("W201").Value = 0.3008144
$ws.Range("X201").Value = 0.3742232
$ws.Range("Y201").Value = 0.2742259
$ws.Range("Z201").Value = 0.6965109
$ws.Range("AA201").Value = 0.5753819
$ws.Range("AB201").Value = 0.9666289
$ws.Range("AD201").Value = 2.2463598
$ws.Range("AE201").Value = 1.1443997
$ws.Range("AF201").Value = 0.881314
$ws.Range("AG201").Value = 0.2699529
$ws.Range("AH201").Value = 0.4071266
$ws.Range("AI201").Value = 0.6104015
$ws.Range("AJ201").Value = 0.0589992
$ws.Range("AK201").Value = 0.0986597
$ws.Range("AL201").Value = 0.428612
$ws.Range("AM201").Value = 0.4041894
$ws.Range("AN201").Value = 0.5924852
$ws.Range("AO201").Value = 0.8304655
$ws.Range("AP201").Value = 0.5044631000000001
$ws.Range("AQ201").Value = 0.4377652
$ws.Range("AS201").Value = 0
$ws.Range("AT201").Value = 0.4306808
$ws.Range("AU201").Value = 1.492323
$ws.Range("AV201").Value = 1.1525576
$ws.Range("AW201").Value = 0.8699282
$ws.Range("AX201").Value = 0.4492076
$ws.Range("AY201").Value = 0.3821432
$ws.Range("BA201").Value = 0
$ws.Range("BB201").Value = 0.225039
$ws.Range("BC201").Value = 0.8363726
$ws.Range("BD201").Value = 0.8251975

# Row 202
$ws.Range("B202").Value = 1.6444436
$ws.Range("C202").Value = 0.8666474
$ws.Range("D202").Value = 2.0388211
$ws.Range("F202").Value = 0.7711825
$ws.Range("G202").Value = 0.6954527
$ws.Range("H202").Value = 0.5419035
$ws.Range("I202").Value = 0.0134373
$ws.Range("K202").Value = 0.446715
$ws.Range("L202").Value = 0.8850055999999999
$ws.Range("M202").Value = 0.5713849
$ws.Range("O202").Value = 0
$ws.Range("P202").Value = 0.5945108
$ws.Range("Q202").Value = 1.2428764
$ws.Range("R202").Value = 0.9103725
$ws.Range("S202").Value = 0.5501423
$ws.Range("T202").Value = 0.6100731
$ws.Range("U202").Value = 1.2701973
$ws.Range("V202").Value = 1.4805952
$ws.Range("W202").Value = 0.1224769
$ws.Range("X202").Value = 1.4282219
$ws.Range("Y202").Value = 0.9415642
$ws.Range("Z202").Value = 0.6273107999999999
$ws.Range("AA202").Value = 0.2515497
$ws.Range("AB202").Value = 0.8272913
$ws.Range("AD202").Value = 0.8060025
$ws.Range("AE202").Value = 0.2657827
$ws.Range("AF202").Value = 0.785845
$ws.Range("AG202").Value = 0.670968
$ws.Range("AH202").Value = 0.9795829
$ws.Range("AI202").Value = 0.3710353
$ws.Range("AJ202").Value = 0.5134248
$ws.Range("AK202").Value = 1.3075544
$ws.Range("AL202").Value = 0.1685297
$ws.Range("AM202").Value = 0.5153522
$ws.Range("AN202").Value = 0.2556862
$ws.Range("AO202").Value = 0.6807695
$ws.Range("AP202").Value = 0.65034
$ws.Range("AQ202").Value = 0.5340679
$ws.Range("AS202").Value = 0.1459652
$ws.Range("AT202").Value = 0.6620194
$ws.Range("AU202").Value = 1.4665934
$ws.Range("AV202").Value = 0.8573318
$ws.Range("AW202").Value = 0.892961
$ws.Range("AX202").Value = 1.1025484
$ws.Range("AY202").Value = 0.7347849
$ws.Range("BA202").Value = 1.1168657
$ws.Range("BB202").Value = 0.5618668999999999
$ws.Range("BC202").Value = 0.5078756
$ws.Range("BD202").Value = 0.5071868

# Row 203
$ws.Range("C203").Value = 0.9722099
$ws.Range("D203").Value = 0.4819602
$ws.Range("F203").Value = 1.4197108
$ws.Range("G203").Value = 0.6840904
$ws.Range("H203").Value = 0.4218212
$ws.Range("I203").Value = 0.7792876
$ws.Range("K203").Value = 0.1836589
$ws.Range("L203").Value = 0.7361999
$ws.Range("M203").Value = 1.3947214
$ws.Range("P203").Value = 1.3444232
$ws.Range("Q203").Value = 1.0852117
$ws.Range("R203").Value = 0.4758084
$ws.Range("S203").Value = 0.4796771
$ws.Range("T203").Value = 0.9927484
$ws.Range("U203").Value = 0.5022966
$ws.Range("V203").Value = 1.4414779
$ws.Range("W203").Value = 0.6391912
$ws.Range("X203").Value = 0.4609018
$ws.Range("Y203").Value = 0
$ws.Range("Z203").Value = 0.4807134
$ws.Range("AA203").Value = 0.7724172
$ws.Range("AB203").Value = 0.4079538
$ws.Range("AD203").Value = 1.0878014
$ws.Range("AE203").Value = 0.8367314
$ws.Range("AF203").Value = 0.3995724
$ws.Range("AH203").Value = 1.3578881
$ws.Range("AI203").Value = 0.7649458
$ws.Range("AJ203").Value = 0.516332
$ws.Range("AK203").Value = 0.2659135
$ws.Range("AL203").Value = 0.2621267
$ws.Range("AM203").Value = 0.4871597
$ws.Range("AN203").Value = 0.7253489
$ws.Range("AO203").Value = 1.2326255
$ws.Range("AP203").Value = 0.4923392
$ws.Range("AQ203").Value = 0.4403353
$ws.Range("AS203").Value = 0
$ws.Range("AT203").Value = 1.3021196
$ws.Range("AV203").Value = 0.831212
$ws.Range("AW203").Value = 0.8989991000000001
$ws.Range("AX203").Value = 0.2920329
$ws.Range("AY203").Value = 0.4765521
$ws.Range("BB203").Value = 0.3874864
$ws.Range("BC203").Value = 0.3940101
$ws.Range("BD203").Value = 0.5800423

# Row 204
$ws.Range("C204").Value = 1.2696023
$ws.Range("D204").Value = 1.2727605
$ws.Range("F204").Value = 0.4853189
$ws.Range("G204").Value = 0.5928349000000001
$ws.Range("H204").Value = 0.7860533
$ws.Range("I204").Value = 0.3886075
$ws.Range("L204").Value = 0.6892722
$ws.Range("M204").Value = 0.9353999
$ws.Range("P204").Value = 0.961294
$ws.Range("Q204").Value = 0.5891006
$ws.Range("R204").Value = 0.9587002
$ws.Range("S204").Value = 0.9480836
$ws.Range("T204").Value = 0.1618068
$ws.Range("U204").Value = 0.1623119
$ws.Range("V204").Value = 1.1846302
$ws.Range("W204").Value = 0.1893046
$ws.Range("X204").Value = 0.4679142
$ws.Range("Y204").Value = 0.7730029
$ws.Range("Z204").Value = 0.8451597
$ws.Range("AA204").Value = 0.0924878
$ws.Range("AB204").Value = 0.3374952
$ws.Range("AD204").Value = 0.5593051999999999
$ws.Range("AF204").Value = 0.4503434
$ws.Range("AH204").Value = 1.1168648
$ws.Range("AI204").Value = 0
$ws.Range("AJ204").Value = 0.1748484
$ws.Range("AK204").Value = 0.7178986000000001
$ws.Range("AL204").Value = 0.5297845
$ws.Range("AM204").Value = 0.4130573
$ws.Range("AN204").Value = 0.6983992999999999
$ws.Range("AO204").Value = 0.8507362000000001
$ws.Range("AP204").Value = 0.6206919
$ws.Range("AQ204").Value = 0.2667475
$ws.Range("AT204").Value = 0.544049
$ws.Range("AV204").Value = 0.9022331
$ws.Range("AW204").Value = 1.1517815
$ws.Range("AX204").Value = 0.7111204
$ws.Range("AY204").Value = 0.8904166999999999
$ws.Range("BB204").Value = 0.0472165
$ws.Range("BC204").Value = 0.9011869
$ws.Range("BD204").Value = 0.3443329

# Row 205
$ws.Range("C205").Value = 0.9796195
$ws.Range("F205").Value = 0.9275067
$ws.Range("G205").Value = 0.5598777
$ws.Range("H205").Value = 0.4974086
$ws.Range("L205").Value = 0.5603941
$ws.Range("M205").Value = 0.6927171
$ws.Range("R205").Value = 0.5671807
$ws.Range("S205").Value = 0.7831838
$ws.Range("U205").Value = 0.7774109
$ws.Range("V205").Value = 0.7642231
$ws.Range("W205").Value = 0.1765279
$ws.Range("X205").Value = 0.4216999
$ws.Range("Z205").Value = 0.4409167
$ws.Range("AA205").Value = 0.6367789
$ws.Range("AB205").Value = 0.2608101
$ws.Range("AF205").Value = 1.0228674
$ws.Range("AJ205").Value = 0.5319709
$ws.Range("AL205").Value = 0.9429176
$ws.Range("AM205").Value = 0.5663383
$ws.Range("AN205").Value = 0.5911761
$ws.Range("AO205").Value = 1.2057493
$ws.Range("AP205").Value = 0.3623678
$ws.Range("AQ205").Value = 0.8477501
$ws.Range("AT205").Value = 0.7187372
$ws.Range("AV205").Value = 0.6205762
$ws.Range("AW205").Value = 0.7348827999999999
$ws.Range("AX205").Value = 0.4511846
$ws.Range("AY205").Value = 0.4069178
$ws.Range("BB205").Value = 0.2773864
$ws.Range("BC205").Value = 0.931844
